$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.670.25'
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').Value = '2.451.37'
$ws.Range('E3').Value = '  -1.65%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.42'
$ws.Range('E5').Value = '  -2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '160.95'
$ws.Range('E6').Value = '  -2.56%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('D9').Value = '2.450.38'
$ws.Range('E9').Value = '  -1.66%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.149'
$ws.Range('E10').Value = '  -7.12%  '
$ws.Range('E11').Value = '  -0.95%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.334'
$ws.Range('E12').Value = '  -5.05%  '
$ws.Range('E13').Value = '  -1.93%  '
$ws.Range('D14').Value = '2.891.44'
$ws.Range('E14').Value = '  -1.90%  '
$ws.Range('D15').Value = '68.509.52'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('E16').Value = '  -4.49%  '
$ws.Range('E17').Value = '  -3.76%  '
$ws.Range('D18').Value = '2.469.26'
$ws.Range('E18').Value = '  -0.45%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.71'
$ws.Range('E19').Value = '  -4.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.21'
$ws.Range('E20').Value = '  -2.00%  '
$ws.Range('E21').Value = '  -5.62%  '
$ws.Range('E22').Value = '  -2.52%  '
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('E24').Value = '  +0.36%  '
$ws.Range('E25').Value = '  -2.56%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '66.43'
$ws.Range('E26').Value = '  -4.08%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.64'
$ws.Range('E27').Value = '  -5.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.07'
$ws.Range('E30').Value = '  -5.89%  '
$ws.Range('D31').Value = '0.0₃0812'
$ws.Range('E31').Value = '  -6.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.13'
$ws.Range('E32').Value = '  -6.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '434.97'
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.13'
$ws.Range('E35').Value = '  -5.09%  '
$ws.Range('E36').Value = '  -6.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '157.23'
$ws.Range('E37').Value = '  +1.23%  '
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('E40').Value = '  -3.18%  '
$ws.Range('E41').Value = '  -1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.302'
$ws.Range('E42').Value = '  -3.50%  '
$ws.Range('E43').Value = '  -3.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '37.41'
$ws.Range('E44').Value = '  -0.84%  '
$ws.Range('E45').Value = '  -7.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.09'
$ws.Range('E46').Value = '  +1.96%  '
$ws.Range('E47').Value = '  -5.95%  '
$ws.Range('E48').Value = '  -4.98%  '
$ws.Range('E49').Value = '  -2.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0714'
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('E51').Value = '  -4.82%  '
